{"js": "// Replace the text of specific table-cell paragraphs (three-digit \u00f7 one-digit\n// division answers) with their updated values. Each \"before\" string is\n// unique in the document, so we can safely match full paragraph text.\nconst replacements = [\n  [\"358\u00f74=89, 2\", \"755\u00f79=83, 8\"],\n  [\"909\u00f72=454, 1\", \"927\u00f79=103, 0\"],\n  [\"353\u00f75=70, 3\", \"357\u00f75=71, 2\"],\n  [\"870\u00f77=124, 2\", \"680\u00f74=170, 0\"],\n  [\"157\u00f74=39, 1\", \"430\u00f79=47, 7\"],\n  [\"225\u00f74=56, 1\", \"252\u00f73=84, 0\"],\n  [\"372\u00f79=41, 3\", \"903\u00f77=129, 0\"],\n  [\"459\u00f73=153, 0\", \"800\u00f77=114, 2\"],\n  [\"347\u00f75=69, 2\", \"446\u00f75=89, 1\"],\n  [\"851\u00f73=283, 2\", \"919\u00f72=459, 1\"],\n  [\"677\u00f79=75, 2\", \"604\u00f76=100, 4\"],\n  [\"597\u00f78=74, 5\", \"683\u00f72=341, 1\"],\n  [\"699\u00f78=87, 3\", \"949\u00f75=189, 4\"],\n  [\"280\u00f77=40, 0\", \"439\u00f75=87, 4\"],\n  [\"209\u00f73=69, 2\", \"199\u00f74=49, 3\"],\n  [\"874\u00f76=145, 4\", \"647\u00f73=215, 2\"],\n  [\"367\u00f73=122, 1\", \"317\u00f74=79, 1\"],\n  [\"747\u00f72=373, 1\", \"149\u00f78=18, 5\"],\n  [\"771\u00f75=154, 1\", \"502\u00f79=55, 7\"],\n  [\"804\u00f77=114, 6\", \"272\u00f76=45, 2\"],\n  [\"139\u00f79=15, 4\", \"381\u00f77=54, 3\"],\n  [\"986\u00f73=328, 2\", \"720\u00f77=102, 6\"],\n  [\"711\u00f77=101, 4\", \"673\u00f75=134, 3\"],\n  [\"692\u00f76=115, 2\", \"337\u00f76=56, 1\"],\n  [\"548\u00f78=68, 4\", \"165\u00f73=55, 0\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst map = new Map(replacements);\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (map.has(text)) {\n    paragraph.insertText(map.get(text), \"Replace\");\n    map.delete(text);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit / one-digit division answers in the table.\n# Each \"before\" string occurs exactly once in the document, so a simple\n# Find/Replace (ReplaceAll) per pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"358\u00f74=89, 2\", \"755\u00f79=83, 8\"),\n    @(\"909\u00f72=454, 1\", \"927\u00f79=103, 0\"),\n    @(\"353\u00f75=70, 3\", \"357\u00f75=71, 2\"),\n    @(\"870\u00f77=124, 2\", \"680\u00f74=170, 0\"),\n    @(\"157\u00f74=39, 1\", \"430\u00f79=47, 7\"),\n    @(\"225\u00f74=56, 1\", \"252\u00f73=84, 0\"),\n    @(\"372\u00f79=41, 3\", \"903\u00f77=129, 0\"),\n    @(\"459\u00f73=153, 0\", \"800\u00f77=114, 2\"),\n    @(\"347\u00f75=69, 2\", \"446\u00f75=89, 1\"),\n    @(\"851\u00f73=283, 2\", \"919\u00f72=459, 1\"),\n    @(\"677\u00f79=75, 2\", \"604\u00f76=100, 4\"),\n    @(\"597\u00f78=74, 5\", \"683\u00f72=341, 1\"),\n    @(\"699\u00f78=87, 3\", \"949\u00f75=189, 4\"),\n    @(\"280\u00f77=40, 0\", \"439\u00f75=87, 4\"),\n    @(\"209\u00f73=69, 2\", \"199\u00f74=49, 3\"),\n    @(\"874\u00f76=145, 4\", \"647\u00f73=215, 2\"),\n    @(\"367\u00f73=122, 1\", \"317\u00f74=79, 1\"),\n    @(\"747\u00f72=373, 1\", \"149\u00f78=18, 5\"),\n    @(\"771\u00f75=154, 1\", \"502\u00f79=55, 7\"),\n    @(\"804\u00f77=114, 6\", \"272\u00f76=45, 2\"),\n    @(\"139\u00f79=15, 4\", \"381\u00f77=54, 3\"),\n    @(\"986\u00f73=328, 2\", \"720\u00f77=102, 6\"),\n    @(\"711\u00f77=101, 4\", \"673\u00f75=134, 3\"),\n    @(\"692\u00f76=115, 2\", \"337\u00f76=56, 1\"),\n    @(\"548\u00f78=68, 4\", \"165\u00f73=55, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
